$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 14773553.80676487
